# Recovered position scoring for extractor.
# Updates the "Status" column (column C) on the "Raw" worksheet for a batch of
# rows that had been left as "Unprocessed" after the extractor ran. Most rows
# are now "Processed"; a couple of rows were found to have been mis-scored and
# are set to "Incorrect" / "Partial" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw")

# Rows that moved from "Unprocessed" -> "Processed"
$processedRows = @(92,93,96,97,100,101,102,103,817,819,820,821,822,823,824,825,826,827,828,829,830,831,832,834,835,837,839,842,850,851,852,853,854,855,861,862,863,864,865,866,867,868,869,870,871,872,873,874,875,876,877,878,879,880,881,882,883,884,885,886,887,888,889,890,891,892,893,894,896,897,898,900,901,902,903,904,908,909,910,911,912,913,914,915,916,917,918,919,920,921,922,923,924,925,927,928,929,930,931,932,933,934,935,936,937,938,939,940,941,942,943,944,945,946,947,948,949,950)

# Rows that moved from "Unprocessed" -> "Incorrect"
$incorrectRows = @(895)

# Rows that moved from "Unprocessed" -> "Partial"
$partialRows = @(899)

foreach ($r in $processedRows) {
    $ws.Range("C$r").Value = "Processed"
}

foreach ($r in $incorrectRows) {
    $ws.Range("C$r").Value = "Incorrect"
}

foreach ($r in $partialRows) {
    $ws.Range("C$r").Value = "Partial"
}
